$d = $word.ActiveDocument

# Replace the case/document number "1314/1245" -> "4214/4213" (single occurrence)
$d.Content.Find.Execute("1314/1245", $true, $false, $false, $false, $false,
                         $true, 1, $false, "4214/4213", 2)

# Replace every date "28.10.2025" -> "29.10.2025" (all occurrences throughout the document)
$d.Content.Find.Execute("28.10.2025", $true, $false, $false, $false, $false,
                         $true, 1, $false, "29.10.2025", 2)
